# Auto update: 2025-11-29 02:59:41
# Rebuilds the scoring table: renames/extends headers, refreshes the
# per-ticker rows with the latest run's numbers, reorders MetLife/AIG,
# drops the stale RSI column, and appends the new macro-prediction columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (row 1) ----------------------------------------------
$ws.Range("G1").Value = "점수(룰)"
$ws.Range("H1").Value = "3일상승확률(%)"
$ws.Range("I1").Value = "5일상승확률(%)"
$ws.Range("J1").Value = "10일상승확률(%)"
$ws.Range("K1").Value = "최종점수"

# New trailing headers - copy the existing header formatting (bold,
# bordered, centered) from K1 onto L1:O1 before writing their text.
$ws.Range("K1").Copy()
$ws.Range("L1:O1").PasteSpecial(-4122)
$ws.Range("L1").Value = "예측방식"
$ws.Range("M1").Value = "판단"
$ws.Range("N1").Value = "MACRO_SCORE"
$ws.Range("O1").Value = "MACRO_SIGNAL"

# ---- Data rows (rows 2-5) ---------------------------------------------
# Column A holds the run date as plain text ("2025-11-29"), not an Excel
# date serial. A direct .Value assignment would get auto-parsed into a
# date serial number, so instead write it as a text-returning formula and
# immediately flatten that formula down to its literal string result -
# this keeps the cell a plain shared-string with no style/number-format
# side effects (matching the original unstyled data cells).
$ws.Range("A2").Formula = '="2025-11-29"'
$ws.Range("A3").Formula = '="2025-11-29"'
$ws.Range("A4").Formula = '="2025-11-29"'
$ws.Range("A5").Formula = '="2025-11-29"'
$ws.Range("A2:A5").Copy()
$ws.Range("A2:A5").PasteSpecial(-4163)  # xlPasteValues

# Row 2: Prudential Financial / PRU
$ws.Range("B2").Value = "Prudential Financial, Inc."
$ws.Range("C2").Value = "PRU"
$ws.Range("D2").Value = 108.24
$ws.Range("E2").ClearContents()
$ws.Range("F2").Value = 3.95
$ws.Range("G2").Value = 50
$ws.Range("H2").Value = 53
$ws.Range("I2").Value = 46
$ws.Range("J2").Value = 40
$ws.Range("K2").Value = 59
$ws.Range("L2").Value = "Pattern"
$ws.Range("M2").Value = "⛔ 관망하십시오."
$ws.Range("N2").Value = 85.36763896678245
$ws.Range("O2").Value = "🟢 완화적 (상승 우위)"

# Row 3: UnitedHealth Group / UNH
$ws.Range("B3").Value = "UnitedHealth Group Incorporated"
$ws.Range("C3").Value = "UNH"
$ws.Range("D3").Value = 330.1
$ws.Range("E3").ClearContents()
$ws.Range("F3").Value = 5.96
$ws.Range("G3").Value = 40
$ws.Range("H3").Value = 46
$ws.Range("I3").Value = 46
$ws.Range("J3").Value = 33
$ws.Range("K3").Value = 56
$ws.Range("L3").Value = "Pattern"
$ws.Range("M3").Value = "⛔ 관망하십시오."
$ws.Range("N3").Value = 85.36763896678245
$ws.Range("O3").Value = "🟢 완화적 (상승 우위)"

# Row 4: MetLife / MET (swapped ahead of AIG)
$ws.Range("B4").Value = "MetLife, Inc."
$ws.Range("C4").Value = "MET"
$ws.Range("D4").Value = 76.58
$ws.Range("E4").ClearContents()
$ws.Range("F4").Value = 3.6
$ws.Range("G4").Value = 30
$ws.Range("H4").Value = 36
$ws.Range("I4").Value = 50
$ws.Range("J4").Value = 36
$ws.Range("K4").Value = 54.6
$ws.Range("L4").Value = "Pattern"
$ws.Range("M4").Value = "⛔ 관망하십시오."
$ws.Range("N4").Value = 85.36763896678245
$ws.Range("O4").Value = "🟢 완화적 (상승 우위)"

# Row 5: American International Group / AIG
$ws.Range("B5").Value = "American International Group, I"
$ws.Range("C5").Value = "AIG"
$ws.Range("D5").Value = 76.18000000000001
$ws.Range("E5").ClearContents()
$ws.Range("F5").Value = 0.64
$ws.Range("G5").Value = 20
$ws.Range("H5").Value = 43
$ws.Range("I5").Value = 50
$ws.Range("J5").Value = 60
$ws.Range("K5").Value = 51.6
$ws.Range("L5").Value = "Pattern"
$ws.Range("M5").Value = "⛔ 관망하십시오."
$ws.Range("N5").Value = 85.36763896678245
$ws.Range("O5").Value = "🟢 완화적 (상승 우위)"
